# feat: add 2022-Q3 data
#
# Target layout:
#   Sheet 1: "总计"    -- summary sheet, gains a new "2022-Q3" row (inserted above
#                          the existing "2022-Q2" row)
#   Sheet 2: "2022-Q3" -- brand new sheet (inserted right before the existing
#                          "2022-Q2" sheet) holding the Q3 fund-holding table
#   Sheet 3: "2022-Q2" -- the original sheet, left untouched

$wb = $excel.ActiveWorkbook

$wsTotal = $wb.Worksheets.Item(1)
$wsQ2    = $wb.Worksheets.Item(2)

# ---------------------------------------------------------------------------
# 1) "总计" sheet: shift the existing 2022-Q2 summary row down to row 3 and
#    insert a new 2022-Q3 summary row at row 2 (2022-Q3 is the newer / first
#    listed quarter).
# ---------------------------------------------------------------------------
$wsTotal.Range("B3").Value = $wsTotal.Range("B2").Value2
$wsTotal.Range("C3").Value = $wsTotal.Range("C2").Value2
$wsTotal.Range("D3").Value = $wsTotal.Range("D2").Value2

# Row index column (0-based, mirrors a pandas DataFrame index): row 2 is
# index 0, the newly-appended row 3 is index 1.
$wsTotal.Range("A3").Value = 1
$wsTotal.Range("A3").Font.Bold = $true
$wsTotal.Range("A3").HorizontalAlignment = -4108
$wsTotal.Range("A3").VerticalAlignment = -4160
$wsTotal.Range("A3").Borders.LineStyle = 1

$wsTotal.Range("A2").Value = 0
$wsTotal.Range("B2").Value = "2022-Q3"
$wsTotal.Range("C2").Value = 18
$wsTotal.Range("D2").Value = 5.77

# ---------------------------------------------------------------------------
# 2) Insert the brand new "2022-Q3" sheet right before the existing
#    "2022-Q2" sheet (which keeps its original data untouched).
# ---------------------------------------------------------------------------
$wsQ3 = $wb.Worksheets.Add($wsQ2)
$wsQ3.Name = "2022-Q3"

$headers = @("基金代码","基金名称","基金规模","股票总仓位","仓位占比","持有市值(亿元)","仓位排名")
for ($i = 0; $i -lt $headers.Length; $i++) {
    $cell = $wsQ3.Cells.Item(1, $i + 2)
    $cell.Value = $headers[$i]
    $cell.Font.Bold = $true
    $cell.HorizontalAlignment = -4108
    $cell.VerticalAlignment = -4160
    $cell.Borders.LineStyle = 1
}

$rows = @(
    @(0,  "010108", "景顺长城核心招景混合A",                    "48.11", "83.11", "4.51", "2.1698", 3),
    @(1,  "010027", "景顺长城核心中景一年持有期混合",            "44.90", "86.84", "4.83", "2.1687", 3),
    @(2,  "009190", "景顺长城核心优选一年持有期混合",            "11.52", "86.01", "8.10", "0.9331", 3),
    @(3,  "159636", "工银瑞信国证港股通科技ETF",                 "8.42",  "98.06", "3.33", "0.2804", 10),
    @(4,  "501021", "华宝标普香港上市中国中小盘指数（LOF）A",    "4.19",  "92.99", "1.94", "0.0813", 6),
    @(5,  "010783", "德邦沪港深龙头混合A",                       "0.55",  "84.96", "5.86", "0.0322", 2),
    @(6,  "013897", "德邦港股通成长精选混合型证券投资基金A",     "0.41",  "79.99", "6.00", "0.0246", 3),
    @(7,  "013898", "德邦港股通成长精选混合型证券投资基金C",     "0.37",  "79.99", "6.00", "0.0222", 3),
    @(8,  "010784", "德邦沪港深龙头混合C",                       "0.36",  "84.96", "5.86", "0.0211", 2),
    @(9,  "519601", "海富通中国海外精选混合（QDII）",            "0.51",  "73.52", "3.01", "0.0154", 8),
    @(10, "501303", "广发恒生中型股指数（LOF）A",                "0.21",  "89.12", "2.26", "0.0047", 3),
    @(11, "006127", "华宝标普香港上市中国中小盘指数（LOF）C",    "0.24",  "92.99", "1.94", "0.0047", 6),
    @(12, "519602", "海富通大中华精选混合（QDII）",              "0.10",  "87.37", "3.90", "0.0039", 8),
    @(13, "003279", "融通沪港深智慧生活灵活配置混合",            "0.10",  "55.71", "2.81", "0.0028", 10),
    @(14, "004996", "广发恒生中型股指数（LOF）C",                "0.09",  "89.12", "2.26", "0.0020", 3),
    @(15, "160922", "大成恒生综合中小型股指数（QDII-LOF）A",     "0.09",  "86.62", "1.71", "0.0015", 3),
    @(16, "015752", "景顺长城核心招景混合C",                     "0.01",  "83.11", "4.51", "0.0005", 3),
    @(17, "008972", "大成恒生综合中小型股指数C",                 "0.02",  "86.62", "1.71", "0.0003", 3)
)

foreach ($r in $rows) {
    $rowIndex = [int]$r[0] + 2

    $cellA = $wsQ3.Cells.Item($rowIndex, 1)
    $cellA.Value = [int]$r[0]
    $cellA.Font.Bold = $true
    $cellA.HorizontalAlignment = -4108
    $cellA.VerticalAlignment = -4160
    $cellA.Borders.LineStyle = 1

    # Fund code (e.g. "010108") must stay text -- a leading apostrophe stops
    # the leading zero from being parsed away as a number.
    $wsQ3.Cells.Item($rowIndex, 2).Value = "'" + $r[1]
    # Fund name is plain text, never numeric-looking.
    $wsQ3.Cells.Item($rowIndex, 3).Value = $r[2]
    # Scale / position / ratio / value columns are stored as TEXT (not
    # numbers) in the source data, so force them with the same trick.
    $wsQ3.Cells.Item($rowIndex, 4).Value = "'" + $r[3]
    $wsQ3.Cells.Item($rowIndex, 5).Value = "'" + $r[4]
    $wsQ3.Cells.Item($rowIndex, 6).Value = "'" + $r[5]
    $wsQ3.Cells.Item($rowIndex, 7).Value = "'" + $r[6]
    # Rank column is a genuine number.
    $wsQ3.Cells.Item($rowIndex, 8).Value = [int]$r[7]
}

Write-Host "2022-Q3 sheet created with $($rows.Length) data rows"
